# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) on Sheet1 held values like "5-9-2012-13";
# correct these to ISO-formatted "2013-05-09" for every data row,
# keeping the cells as plain text (not auto-converted to a date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2013-05-09"

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    # Write the corrected value through a text-literal formula first so
    # Excel's autodetect doesn't coerce the ISO-looking string into a date
    # serial number, then paste-special just the value back onto itself so
    # the formula collapses to a plain text constant (matching the original
    # cell's plain string storage / formatting).
    $cell.Formula = "=""" + $correctDate + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
